# Refresh the cryptos price/volume snapshot (and a few rank swaps) to match
# the latest scrape, mirroring the GitHub Actions commit's OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.073.90'
$ws.Range("E2").Value = '  -0.84%  '
# Row 3
$ws.Range("D3").Value = '1.650.55'
$ws.Range("E3").Value = '  -0.82%  '
# Row 4
$ws.Range("E4").Value = '  -0.41%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.47'
$ws.Range("E5").Value = '  -0.66%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5216'
$ws.Range("E6").Value = '  -2.56%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.006'
$ws.Range("E7").Value = '  -0.36%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2616'
$ws.Range("E8").Value = '  -1.55%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06284'
$ws.Range("E9").Value = '  -1.88%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.53'
$ws.Range("E10").Value = '  -0.30%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07804'
$ws.Range("E11").Value = '  -0.33%  '
# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.682.30'
$ws.Range("E12").Value = '  +0.68%  '
# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.476'
$ws.Range("E13").Value = '  -1.97%  '
# Row 14
$ws.Range("D14").Value = '1.877.11'
$ws.Range("E14").Value = '  -0.78%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5547'
$ws.Range("E15").Value = '  +0.38%  '
# Row 16
$ws.Range("D16").Value = '0.0₅8003'
$ws.Range("E16").Value = '  -2.47%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.74'
$ws.Range("E17").Value = '  -1.52%  '
# Row 18
$ws.Range("D18").Value = '26.063.69'
$ws.Range("E18").Value = '  -0.90%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.633'
$ws.Range("E20").Value = '  -0.99%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '194.25'
$ws.Range("E21").Value = '  +0.14%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.07'
$ws.Range("E22").Value = '  -1.44%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.947'
$ws.Range("E23").Value = '  -1.58%  '
# Row 24
$ws.Range("E24").Value = '  -0.41%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.70'
$ws.Range("E25").Value = '  +0.36%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1203'
$ws.Range("E26").Value = '  -2.35%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.177'
$ws.Range("E27").Value = '  -0.21%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.88'
$ws.Range("E28").Value = '  -1.39%  '
# Row 29
$ws.Range("E29").Value = '  -0.47%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05606'
$ws.Range("E30").Value = '  -4.13%  '
# Row 31
$ws.Range("E31").Value = '  -1.71%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.486'
$ws.Range("E32").Value = '  -3.28%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.356'
$ws.Range("E33").Value = '  +2.23%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.592'
$ws.Range("E34").Value = '  -1.30%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.797'
$ws.Range("E35").Value = '  -1.04%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9485'
$ws.Range("E36").Value = '  -1.62%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.407'
$ws.Range("E37").Value = '  -0.54%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5649'
$ws.Range("E38").Value = '  -2.76%  '
# Row 39
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.966'
$ws.Range("E39").Value = '  +1.56%  '
# Row 40
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01582'
$ws.Range("E40").Value = '  -1.54%  '
# Row 41
$ws.Range("D41").Value = '1.058.77'
$ws.Range("E41").Value = '  +0.63%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.005'
$ws.Range("E42").Value = '  -0.46%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8376'
$ws.Range("E43").Value = '  -3.38%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.99'
$ws.Range("E44").Value = '  -2.64%  '
# Row 45
$ws.Range("D45").Value = '1.788.12'
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.13'
$ws.Range("E46").Value = '  -1.34%  '
# Row 47
$ws.Range("E47").Value = '  -0.74%  '
# Row 48
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05365'
$ws.Range("E48").Value = '  +3.87%  '
# Row 49
$ws.Range("B49").Value = 'Frax'
$ws.Range("C49").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.005'
$ws.Range("E49").Value = '  -0.95%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4330'
$ws.Range("E50").Value = '  -1.24%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.921'
$ws.Range("E51").Value = '  -1.11%  '
